# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets to reflect the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 7601
$wsExhibit.Range("F4").Value = 217
$wsExhibit.Range("F5").Value = 20
$wsExhibit.Range("F6").Value = 264
$wsExhibit.Range("F7").Value = 1143
$wsExhibit.Range("F8").Value = 200
$wsExhibit.Range("F10").Value = 144
$wsExhibit.Range("F11").Value = 37

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7601
$wsAll.Range("F4").Value = 217
$wsAll.Range("F5").Value = 20
$wsAll.Range("F6").Value = 264
$wsAll.Range("F7").Value = 1143
$wsAll.Range("F8").Value = 200
$wsAll.Range("F11").Value = 144
$wsAll.Range("F12").Value = 37
